# Atualizei dados da bibi
# Insert a new daily record (Dia 16, 28514.4, Mes 6, Ano 2025, Periodo 06/2025)
# into the faturamento diario sheet, right after the existing "Dia 15" row
# for June/2025 (row 16), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 17, shifting rows 17..108 down to 18..109
$ws.Rows.Item(17).Insert()

# Populate the new row with the new data point
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 28514.4
$ws.Cells.Item(17, 3).Value = 6
$ws.Cells.Item(17, 4).Value = 2025
$ws.Cells.Item(17, 5).Value = "06/2025"
